$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying data rows for "NOAA's National Weather Service", "Tuckahoe Train
# Station The Flood of Apr 16 2007", "More than 100,000 without power as storm hits
# Ontario, Quebec" and "Patriot's Day Storm Packs a 156-mph Punch" (together with
# their matching uri column) get re-ordered (as part of reprocessing the source
# JSON records for the time-bucket analysis). Re-assign the title (column A) and
# uri (column E) text for rows 3-6 so each title stays paired with its own uri,
# but in the new row order: Tuckahoe, Patriot's Day, NOAA's, More than 100,000.

$ws.Cells.Item(3, 1).Value2 = "Tuckahoe Train Station The Flood of Apr 16 2007"
$ws.Cells.Item(3, 5).Value2 = "https://www.flickr.com/photos/kc2gog/sets/72157604303045368/with/2372814086/"

$ws.Cells.Item(4, 1).Value2 = "Patriot's Day Storm Packs a 156-mph Punch"
$ws.Cells.Item(4, 5).Value2 = "http://www.mountwashington.org/news/release.php?id=17"

$ws.Cells.Item(5, 1).Value2 = "NOAA's National Weather Service"
$ws.Cells.Item(5, 5).Value2 = "http://www.erh.noaa.gov/er/bgm/WeatherEvents/Snow/april162007/april162007.shtml"

$ws.Cells.Item(6, 1).Value2 = "More than 100,000 without power as storm hits Ontario, Quebec"
$ws.Cells.Item(6, 5).Value2 = "http://www.cbc.ca/canada/story/2007/04/16/spring-storm-mon.html"
